# Timesheet now signed off by the supervisor:
#  - Supervisor Name filled in (G6)
#  - Supervisor Signature initials + sign-off date added (A27 / D27)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor Signature (initials) and Date, mirroring the employee
# sign-off row (A25/D25) directly above the "Supervisor Signature" label.
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = 41682
$ws.Range("D27").NumberFormat = "mm-dd-yy"

# Leave the selection where the author left it after signing off.
[void]$ws.Range("F34").Select()
